# Update the "想去人数" (want-to-go count) figures in column F
# on the "展览" and "全部类型" worksheets, reflecting refreshed
# data from the site scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 4499
$wsExhibit.Range("F5").Value = 50
$wsExhibit.Range("F9").Value = 2619
$wsExhibit.Range("F11").Value = 214
$wsExhibit.Range("F12").Value = 82
$wsExhibit.Range("F13").Value = 5359
$wsExhibit.Range("F15").Value = 219
$wsExhibit.Range("F16").Value = 567
$wsExhibit.Range("F17").Value = 11510
$wsExhibit.Range("F18").Value = 11625

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 4499
$wsAll.Range("F5").Value = 50
$wsAll.Range("F9").Value = 2619
$wsAll.Range("F12").Value = 214
$wsAll.Range("F13").Value = 82
$wsAll.Range("F14").Value = 5359
$wsAll.Range("F16").Value = 219
$wsAll.Range("F17").Value = 567
$wsAll.Range("F18").Value = 11510
$wsAll.Range("F19").Value = 11625
